$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count() + $usedRange.Row() - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $value = $cell.Value()
    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ",\s*"

        $hasSystem = $false
        foreach ($p in $parts) {
            if ($p.Equals("System")) { $hasSystem = $true }
        }

        if ($hasSystem -and -not $parts[0].Equals("System")) {
            $rest = @()
            foreach ($p in $parts) {
                if (-not $p.Equals("System")) { $rest += $p }
            }
            $newParts = @("System") + $rest
            $cell.Value = [string]::Join(", ", $newParts)
        }
    }
}
